$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "HDBSCAN vs DBSCAN"
$ws.Range("E32").Value = "https://dodonam.tistory.com/326"

$ws.Range("D36").Value = "타이어 산업 데이터 특징 및 성능 예측 사례"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/324"

$ws.Range("D50").Value = "주피터허브?"
$ws.Range("E50").Value = "http://incredible.egloos.com/7521342"
